$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 4, 3, 16),
    @(6, 8, 7, 12),
    @(5, 7, 4, 13),
    @(3, 5, 4, 15),
    @(2, 15, 3, 5),
    @(5, 4, 3, 16),
    @(5, 8, 1, 12),
    @(5, 13, 9, 7),
    @(2, 8, 3, 12),
    @(3, 12, 4, 8),
    @(3, 17, 4, 3),
    @(4, 6, 6, 14)
)

$startRow = 1035
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}

$ws.Range("A1047").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1028
